# Insert a new weekly record at the top of the "Macroferia Regional de Talca -
# Pepino ensalada" data block (row 644), pushing the existing rows 644-668
# down to 645-669.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 644:668 down by one row to make room for the new record.
$ws.Rows.Item(644).Insert()

# Populate the newly inserted row 644 with the new weekly entry.
$ws.Range("A644").Value = 5
$ws.Range("B644").Value = "Macroferia Regional de Talca"
$ws.Range("C644").Value = "Maule"
$ws.Range("D644").Value = 45075
$ws.Range("E644").Value = 7
$ws.Range("F644").Value = 100112043
$ws.Range("G644").Value = "Pepino ensalada"
$ws.Range("H644").Value = "Sin especificar"
$ws.Range("I644").Value = "Primera"
$ws.Range("J644").Value = 500
$ws.Range("K644").Value = 10000
$ws.Range("L644").Value = 10000
$ws.Range("M644").Value = 10000
$ws.Range("N644").Value = "$/caja 60 unidades"
$ws.Range("O644").Value = "Región de Arica y Parinacota"
$ws.Range("P644").Value = 167
$ws.Range("Q644").Value = 60
$ws.Range("R644").Value = "Hortaliza"
